$d = $word.ActiveDocument

# The title paragraph ("Answers Gaussian elimination") was authored as five
# separate runs: "Answers" / " " / "Gaussian" / " " / "elimination".
# Collapse them into a single run containing the full text, leaving every
# other paragraph (and the bookmark around it) untouched.
$titleRange = $d.Paragraphs(1).Range
$titleRange.Find.Execute("Answers Gaussian elimination", $false, $false, $false, $false, $false, $true, 1, $false, "Answers Gaussian elimination", 2)
